# Slide 2 ("Table of Content" slide):
#  1) Move the "TextBox 22" shape (currently at 786909, 1176298 EMU) to
#     560767, 588149 EMU. Shape.Left/.Top are expressed in points, so the
#     EMU targets below are the nearest points value (as a 32-bit Single,
#     matching the real PowerPoint object model) that round-trips to the
#     exact target EMU.
#  2) Fix the capitalization of "mini-columns" -> "Mini-columns" in the
#     third bullet of the "Content Placeholder 2" text box, editing the
#     run directly so the paragraph keeps a single run/formatting.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$textBox = $s.Shapes.Item(3)
$textBox.Left = 44.154884338378906
$textBox.Top = 46.31094741821289

$contentPlaceholder = $s.Shapes.Item(2)
$run = $contentPlaceholder.TextFrame.TextRange.Paragraphs(3).Runs(1)
$run.Text = "Investigation of Inactive Mini-columns"
